$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing corrected values (rows 120, 125, 127)
$ws.Range("B120").Value = 9848.700000000001
$ws.Range("C120").Value = 9115.6

$ws.Range("B125").Value = 8136.7
$ws.Range("C125").Value = 7142.6

$ws.Range("B127").Value = 8257.200000000001
$ws.Range("C127").Value = 7191.4

# Append new row 138 with data for 01-07-2021
$ws.Range("A138").Value = "'01-07-2021"
$ws.Range("A138").Style = "Normal"
$ws.Range("B138").Value = 8948.299999999999
$ws.Range("C138").Value = 8148.9
$ws.Range("D138").Value = 8.9
